# Update countries & provincias Spain
# - refresh "Datos actualizados..." timestamp (A1)
# - refresh COVID-19 per-country counters (columns B..H) for the countries
#   whose figures moved in this data pull
# - three country-name triples got re-ranked (their stats leapfrogged a
#   neighbouring row), so the "Pais" label (column A) for those rows is
#   rewritten along with the numbers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Septiembre de 2020 a las 23:18"

$ws.Cells.Item(4, 2).Value = 7072808
$ws.Cells.Item(4, 3).Value = 24696
$ws.Cells.Item(4, 4).Value = 4323113
$ws.Cells.Item(4, 5).Value = 2544482
$ws.Cells.Item(4, 7).Value = 711
$ws.Cells.Item(4, 8).Value = 205213

$ws.Cells.Item(5, 2).Value = 5640496
$ws.Cells.Item(5, 3).Value = 80391
$ws.Cells.Item(5, 4).Value = 4581820
$ws.Cells.Item(5, 5).Value = 968655

$ws.Cells.Item(6, 2).Value = 4591604
$ws.Cells.Item(6, 3).Value = 31521
$ws.Cells.Item(6, 5).Value = 566297
$ws.Cells.Item(6, 7).Value = 758
$ws.Cells.Item(6, 8).Value = 138108

$ws.Cells.Item(12, 2).Value = 663282
$ws.Cells.Item(12, 3).Value = 1346
$ws.Cells.Item(12, 4).Value = 592904
$ws.Cells.Item(12, 5).Value = 54260
$ws.Cells.Item(12, 7).Value = 126
$ws.Cells.Item(12, 8).Value = 16118

$ws.Cells.Item(25, 2).Value = 277176
$ws.Cells.Item(25, 3).Value = 1625
$ws.Cells.Item(25, 5).Value = 21385
$ws.Cells.Item(25, 7).Value = 10
$ws.Cells.Item(25, 8).Value = 9491

$ws.Cells.Item(54, 1).Value = "Costa Rica"
$ws.Cells.Item(54, 2).Value = 66689
$ws.Cells.Item(54, 3).Value = 1087
$ws.Cells.Item(54, 4).Value = 25706
$ws.Cells.Item(54, 5).Value = 40223
$ws.Cells.Item(54, 7).Value = 15
$ws.Cells.Item(54, 8).Value = 760

$ws.Cells.Item(55, 1).Value = "Nepal"
$ws.Cells.Item(55, 2).Value = 66632
$ws.Cells.Item(55, 3).Value = 1356
$ws.Cells.Item(55, 4).Value = 48061
$ws.Cells.Item(55, 5).Value = 18142
$ws.Cells.Item(55, 7).Value = 2
$ws.Cells.Item(55, 8).Value = 429

$ws.Cells.Item(56, 1).Value = "Barein"
$ws.Cells.Item(56, 2).Value = 66402
$ws.Cells.Item(56, 3).Value = 650
$ws.Cells.Item(56, 4).Value = 59367
$ws.Cells.Item(56, 5).Value = 6808
$ws.Cells.Item(56, 7).Value = 3
$ws.Cells.Item(56, 8).Value = 227

$ws.Cells.Item(84, 2).Value = 19343
$ws.Cells.Item(84, 3).Value = 16
$ws.Cells.Item(84, 4).Value = 18682
$ws.Cells.Item(84, 5).Value = 541

$ws.Cells.Item(92, 2).Value = 13578
$ws.Cells.Item(92, 3).Value = 23
$ws.Cells.Item(92, 5).Value = 5982

$ws.Cells.Item(100, 2).Value = 9818
$ws.Cells.Item(100, 3).Value = 48
$ws.Cells.Item(100, 4).Value = 8438

$ws.Cells.Item(109, 2).Value = 7403
$ws.Cells.Item(109, 3).Value = 19
$ws.Cells.Item(109, 4).Value = 7011
$ws.Cells.Item(109, 5).Value = 231

$ws.Cells.Item(116, 2).Value = 5407
$ws.Cells.Item(116, 3).Value = 3
$ws.Cells.Item(116, 4).Value = 5338
$ws.Cells.Item(116, 5).Value = 8

$ws.Cells.Item(127, 2).Value = 4738
$ws.Cells.Item(127, 3).Value = 16
$ws.Cells.Item(127, 4).Value = 2991
$ws.Cells.Item(127, 5).Value = 1720

$ws.Cells.Item(129, 2).Value = 4236
$ws.Cells.Item(129, 3).Value = 119
$ws.Cells.Item(129, 4).Value = 1462
$ws.Cells.Item(129, 5).Value = 2619
$ws.Cells.Item(129, 7).Value = 1
$ws.Cells.Item(129, 8).Value = 155

$ws.Cells.Item(132, 1).Value = "Siria"
$ws.Cells.Item(132, 2).Value = 3877
$ws.Cells.Item(132, 3).Value = 44
$ws.Cells.Item(132, 4).Value = 983
$ws.Cells.Item(132, 5).Value = 2716
$ws.Cells.Item(132, 7).Value = 3
$ws.Cells.Item(132, 8).Value = 178

$ws.Cells.Item(133, 1).Value = "Lituania"
$ws.Cells.Item(133, 2).Value = 3859
$ws.Cells.Item(133, 3).Value = 45
$ws.Cells.Item(133, 4).Value = 2225
$ws.Cells.Item(133, 5).Value = 1547
$ws.Cells.Item(133, 8).Value = 87

$ws.Cells.Item(134, 2).Value = 3665
$ws.Cells.Item(134, 3).Value = 78
$ws.Cells.Item(134, 4).Value = 2426
$ws.Cells.Item(134, 5).Value = 1214
$ws.Cells.Item(134, 7).Value = 1
$ws.Cells.Item(134, 8).Value = 25

$ws.Cells.Item(153, 4).Value = 1240
$ws.Cells.Item(153, 5).Value = 202

$ws.Cells.Item(157, 1).Value = "Togo"
$ws.Cells.Item(157, 2).Value = 1683
$ws.Cells.Item(157, 3).Value = 14
$ws.Cells.Item(157, 4).Value = 1290
$ws.Cells.Item(157, 5).Value = 352
$ws.Cells.Item(157, 8).Value = 41

$ws.Cells.Item(158, 1).Value = "Principado de Andorra"
$ws.Cells.Item(158, 2).Value = 1681
$ws.Cells.Item(158, 4).Value = 1199
$ws.Cells.Item(158, 5).Value = 429
$ws.Cells.Item(158, 8).Value = 53

$ws.Cells.Item(189, 2).Value = 210
$ws.Cells.Item(189, 3).Value = 1
$ws.Cells.Item(189, 4).Value = 205

$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1

$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
